# Issue lifecycle diagram: rename the "Accepted" and "Ongoing" (first
# instance) status badges to "Prioritized" and "Assigned" respectively.
# Each badge's text is built from two runs -- a leading "s." run (plain
# weight) and the status-word run (bold). We fold the new word into the
# bold run and clear the leading "s." run so the shape ends up with a
# single bold run holding the new word, matching how the author's edit
# collapsed the two runs into one.

function Set-BadgeText($ShapeId, $NewText) {
    $p = $ppt.ActivePresentation
    $s = $p.Slides.Item(1)

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Id -eq $ShapeId) {
            $tr = $shp.TextFrame.TextRange
            $runCount = $tr.Runs().Count
            if ($runCount -ge 2) {
                # Put the new text in the (bold) second run, then blank
                # out the leading "s." run so it collapses away on save.
                $tr.Runs(2).Text = $NewText
                $tr.Runs(1).Text = ""
            } else {
                $tr.Text = $NewText
            }
            return
        }
    }
}

# Rounded Rectangle 134 (id 135): "s." + "Accepted" -> "Prioritized"
Set-BadgeText 135 "Prioritized"

# Rounded Rectangle 176 (id 177): "s." + "Ongoing" -> "Assigned"
Set-BadgeText 177 "Assigned"
